$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bellarmine
$ws.Range("A2").Value = "Bellarmine"
$ws.Range("B2").Value = "Atlantic Sun"
$ws.Range("C2").Value = "16-13"
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "194"
$ws.Range("D2").ClearFormats()

# Row 3: Iowa State
$ws.Range("A3").Value = "Iowa State"
$ws.Range("B3").Value = "Big 12"
$ws.Range("C3").Value = "20-12"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "47"
$ws.Range("D3").ClearFormats()
